# Normalize the "aba" sheet text:
#  - Column A (provider name) -> UPPERCASE (no-op for rows already upper)
#  - B2/C2 ("fono lar"/"cooperado") -> UPPERCASE
#  - Column D ("Fonoaudiologo a") -> "FONOAUDIOLOGO (A)" for every data row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("aba")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1)
    $aText = $a.Text
    if ($aText) {
        $a.Value = $aText.ToUpper()
    }

    $b = $ws.Cells.Item($r, 2)
    $bText = $b.Text
    if ($bText) {
        $b.Value = $bText.ToUpper()
    }

    $c = $ws.Cells.Item($r, 3)
    $cText = $c.Text
    if ($cText) {
        $c.Value = $cText.ToUpper()
    }

    $d = $ws.Cells.Item($r, 4)
    $dText = $d.Text
    if ($dText -and $dText.ToUpper() -eq "FONOAUDIOLOGO A") {
        $d.Value = "FONOAUDIOLOGO (A)"
    }
}
